# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" block (rows 16-71, columns E and F) for
# the worker JHON JAVIER MEZA QUIROZ is re-sorted from descending period
# order (2102 -> 1607) to ascending period order (1607 -> 2102), carrying
# each period's "Valor Mora" (F) along with it, and the values for the new
# part-1 periods (1809-2102) are refreshed to 31249 (22916 for the most
# recent period, 2102).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1607","1608","1609","1610","1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102")
$values  = @(25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,25774,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,22916)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

Write-Output "Reordered periods/values for rows $startRow..$($startRow + $periods.Length - 1)"
